# Generate Report for Handoff
# Adds two new localization entries (98a512ad-... and f588d7b8-...) as new
# rows 6/7 on the "Overview" sheet and on each per-language detail sheet
# ("zh-cn", "de-de").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# Columns: A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/98a512ad5466460097cbd27eadf23803placeholder/e2e/98a512ad-5466-4600-97cb-d27eadf23803.md", "", "", "98a512ad-5466-4600-97cb-d27eadf23803.md") | Out-Null
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-29-19 12:29:29"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/f588d7b86a0543c0ab6e549d5dea2b01placeholder/e2e/f588d7b8-6a05-43c0-ab6e-549d5dea2b01.md", "", "", "f588d7b8-6a05-43c0-ab6e-549d5dea2b01.md") | Out-Null
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-29-19 12:29:29"

# ---------------------------------------------------------------------
# Helper data shared by both per-language detail sheets
# ---------------------------------------------------------------------
function Add-DetailRows($ws, $lang, $row6Target, $row6TargetDatetime, $row7Target, $row7TargetDatetime) {
    # --- Row 6: 98a512ad-5466-4600-97cb-d27eadf23803 ---
    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/98a512ad5466460097cbd27eadf23803placeholder/e2e/98a512ad-5466-4600-97cb-d27eadf23803.md", "", "", "98a512ad-5466-4600-97cb-d27eadf23803.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/98a512ad5466460097cbd27eadf23803placeholder/e2e/98a512ad-5466-4600-97cb-d27eadf23803.md", "", "", ".md") | Out-Null
    $ws.Range("C6").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D6"), $row6Target, "", "", "98a512ad-5466-4600-97cb-d27eadf23803.2c8cec09f8da8a9f8e021b99081151649c271678.$lang.xlf") | Out-Null
    $ws.Range("E6").Value = $row6TargetDatetime
    $ws.Range("H6").Value = "0001-01-01 00:00:00"
    $ws.Range("I6").Value = "Include"

    # --- Row 7: f588d7b8-6a05-43c0-ab6e-549d5dea2b01 ---
    $ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/f588d7b86a0543c0ab6e549d5dea2b01placeholder/e2e/f588d7b8-6a05-43c0-ab6e-549d5dea2b01.md", "", "", "f588d7b8-6a05-43c0-ab6e-549d5dea2b01.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/f588d7b86a0543c0ab6e549d5dea2b01placeholder/e2e/f588d7b8-6a05-43c0-ab6e-549d5dea2b01.md", "", "", ".md") | Out-Null
    $ws.Range("C7").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D7"), $row7Target, "", "", "f588d7b8-6a05-43c0-ab6e-549d5dea2b01.5afe42b228e43bd8f1fa8fd43f691a2820366f5c.$lang.xlf") | Out-Null
    $ws.Range("E7").Value = $row7TargetDatetime
    $ws.Range("H7").Value = "0001-01-01 00:00:00"
    $ws.Range("I7").Value = "Include"
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Add-DetailRows $wsZhCn "zh-cn" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/98a512ad5466460097cbd27eadf23803placeholder/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/98a512ad-5466-4600-97cb-d27eadf23803.2c8cec09f8da8a9f8e021b99081151649c271678.zh-cn.xlf" `
    "2016-03-19 12:29:26" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f588d7b86a0543c0ab6e549d5dea2b01placeholder/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f588d7b8-6a05-43c0-ab6e-549d5dea2b01.5afe42b228e43bd8f1fa8fd43f691a2820366f5c.zh-cn.xlf" `
    "2016-03-19 12:29:26"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Add-DetailRows $wsDeDe "de-de" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/98a512ad5466460097cbd27eadf23803placeholder/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/98a512ad-5466-4600-97cb-d27eadf23803.2c8cec09f8da8a9f8e021b99081151649c271678.de-de.xlf" `
    "2016-03-19 12:29:29" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f588d7b86a0543c0ab6e549d5dea2b01placeholder/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f588d7b8-6a05-43c0-ab6e-549d5dea2b01.5afe42b228e43bd8f1fa8fd43f691a2820366f5c.de-de.xlf" `
    "2016-03-19 12:29:29"
